$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.498.34"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +3.73%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.499.21"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.32%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.04%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'590.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +3.28%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'169.40"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +5.27%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D8').Value = "'3.497.16"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +2.21%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.587"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +6.26%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'7.32"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.46%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +5.10%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +3.40%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'4.107.14"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.37%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.48%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'28.29"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +4.85%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +2.17%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'66.537.33"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +3.71%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'3.545.54"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +4.48%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'6.32"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +4.21%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +3.94%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'390.32"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.82%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'7.96"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.99%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'73.00"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.23%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.535"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +3.54%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.0000122"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +6.01%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +10.33%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +2.12%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.14%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +5.67%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +6.55%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +2.71%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'23.60"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +3.40%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +4.72%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.00%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +7.26%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'162.36"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +1.75%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +3.54%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'1.90"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +4.95%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'6.86"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +6.82%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'4.69"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +6.28%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'27.56"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +6.38%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +2.92%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'26.46"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +2.82%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'2.801.09"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.11%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'43.13"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.55%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +2.24%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'2.51"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +4.58%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'353.39"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +5.62%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +3.50%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'33.69"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +12.67%  "
$ws.Range('E51').Style = 'Normal'
